$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" on every sheet that uses it ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" status columns ---
# Target stored width is ~13.41 characters; this runtime quantizes ColumnWidth to
# 1/6-character steps (pixel granularity) before it gets serialized back to the
# sheet XML, so 12.5 is the input that lands on the closest achievable width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
